# PARC RZB (version 1).xlsx - "Add files via upload"
#
# Observed changes (from the OOXML diff):
#   1. Six rows on "Feuil1" (D199:D204) have their AGENCE value changed
#      from "EES" to "AQUALIS".
#   2. The sheet selection / active-tab state moved: "Feuil2" was the
#      active/selected sheet (cell H11 selected) before the edit; after the
#      edit "Feuil1" is the active/selected sheet (cell D16 selected), and
#      "Feuil2" is left with cell F82 selected.
#   3. (The pivot cache would normally be refreshed/bumped after the source
#      data edit - that refresh is driven by Excel itself, not by explicit
#      user action, so we just perform the data edit + navigation here.)

$wb = $excel.ActiveWorkbook

$wsParc  = $wb.Worksheets.Item("Feuil2")   # pivot-table sheet
$wsData  = $wb.Worksheets.Item("Feuil1")   # raw data sheet

# --- 1. Update the AGENCE column for the EES -> AQUALIS reclassification ---
$wsData.Range("D199:D204").Value = "AQUALIS"

# Keep the pivot cache/table in sync with the edited source range.
$wb.PivotCaches().Item(1).Refresh()

# --- 2. Leave "Feuil2" with F82 selected (it is no longer the active tab) ---
$wsParc.Select()
$wsParc.Range("F82").Select()

# --- 3. Activate "Feuil1" and select D16 - it becomes the final/active tab ---
$wsData.Select()
$wsData.Range("D16").Select()
